# FTT_variables.xlsx - "Add battery learning, fix gamma values"
#
# On the FTT-Tr sheet, the gamma "switched on/off" flags in column B were
# bugged (fixed at 1) for several rows. Turn those flags off (0) so the
# masterfiles->csv export script actually picks up the FTT-Tr gamma values,
# and update the sheet's selection/scroll state to match.

$wb = $excel.ActiveWorkbook

# The edited sheet is FTT-Tr (the workbook's active sheet / second tab).
$ws = $wb.Worksheets.Item("FTT-Tr")
$ws.Activate()

# Column B holds 0/1 switches next to each variable row. Rows 3-12 and 17
# were incorrectly left "on" (1); flip them to "off" (0).
$rowsToClear = @(3,4,5,6,7,8,9,10,11,12,17)
foreach ($r in $rowsToClear) {
    $ws.Cells.Item($r, 2).Value = 0
}

# Update the sheet view: scroll down a bit and select B4:B22 (anchor cell B4)
# to match the saved selection state in the workbook.
$ws.Range("B4:B22").Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
